# stopwords_Hayatt.xlsx — replace the trailing "니엘 / 호텔 / 롯데" stop-word
# block with the Park Hyatt Busan Downtown terms, extending the list by
# three rows (681 -> 684) and leaving the selection one row past the new
# last entry, just like a user who typed the new words and pressed Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A677").Value = "서울"
$ws.Range("A678").Value = "것"
$ws.Range("A679").Value = "파크"
$ws.Range("A680").Value = "하얏트"
$ws.Range("A681").Value = "부산"
$ws.Range("A682").Value = "다운"
$ws.Range("A683").Value = "파크하얏트"
$ws.Range("A684").Value = "호텔"

$ws.Range("A685").Select()
